# Update "想去人数" (F column) figures across the sheets to reflect the
# newer scrape snapshot used to regenerate the gh-pages output.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F6").Value = 5688
$ws.Range("F7").Value = 527
$ws.Range("F8").Value = 749
$ws.Range("F9").Value = 984
$ws.Range("F10").Value = 830
$ws.Range("F17").Value = 1923
$ws.Range("F18").Value = 1494
$ws.Range("F19").Value = 989
$ws.Range("F23").Value = 583
$ws.Range("F24").Value = 182
$ws.Range("F28").Value = 3183
$ws.Range("F30").Value = 116
$ws.Range("F31").Value = 77
$ws.Range("F34").Value = 435
$ws.Range("F39").Value = 307
$ws.Range("F40").Value = 767
$ws.Range("F41").Value = 98
$ws.Range("F42").Value = 63

# --- Sheet "演出" ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 225

# --- Sheet "全部类型" ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F6").Value = 5688
$ws.Range("F7").Value = 527
$ws.Range("F8").Value = 749
$ws.Range("F10").Value = 225
$ws.Range("F11").Value = 984
$ws.Range("F12").Value = 830
$ws.Range("F22").Value = 1923
$ws.Range("F23").Value = 1494
$ws.Range("F24").Value = 989
$ws.Range("F28").Value = 583
$ws.Range("F29").Value = 182
$ws.Range("F31").Value = 3183
$ws.Range("F33").Value = 116
$ws.Range("F34").Value = 77
$ws.Range("F37").Value = 435
$ws.Range("F41").Value = 307
$ws.Range("F42").Value = 767
$ws.Range("F43").Value = 98
$ws.Range("F44").Value = 63
